$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Current Quotes")

# ---------------------------------------------------------------
# 1. Insert a new column G ("Cost Quoted+Extra Parts/Quantities")
#    Old G (Cost/Unit) -> H, old H (Notes) -> I
# ---------------------------------------------------------------
$ws.Columns.Item(7).Insert()

# ---------------------------------------------------------------
# 2. Insert two new rows before the old row 6 (Osh Park row) so it
#    becomes row 8, creating blank rows 6 & 7 for new quotes.
# ---------------------------------------------------------------
$ws.Rows.Item(6).Resize(2).Insert()

# ---------------------------------------------------------------
# 3. Insert three new rows after row 8 for the panelized gerber
#    quote tiers (rows 9, 10, 11).
# ---------------------------------------------------------------
$ws.Rows.Item(9).Resize(3).Insert()

# ---------------------------------------------------------------
# 4. Numbers / dates / formulas first (order doesn't affect the
#    shared-string table). Text values are written further below
#    in a specific sequence that matches the original authoring
#    order of the workbook.
# ---------------------------------------------------------------

# Existing quote rows 3-5: "Extra Parts" (G) and "Cost/Unit" (H)
$ws.Cells.Item(3, 7).Formula = "=F3+89.28+G7"
$ws.Cells.Item(4, 7).Formula = "=F4+148.8+G7*2"
$ws.Cells.Item(5, 7).Formula = "=F5+89.28+F8*4"
$ws.Cells.Item(3, 8).Formula = "=G3/E3"
$ws.Cells.Item(4, 8).Formula = "=G4/E4"
$ws.Cells.Item(5, 8).Formula = "=G5/E5"

# Row 6: Screaming Circuits panelized assembly quote
$ws.Cells.Item(6, 2).Value = "8/29/2016"
$ws.Cells.Item(6, 5).Value = 12
$ws.Cells.Item(6, 6).Value = 2331.33
$ws.Cells.Item(6, 7).Formula = "=F6+90.28"
$ws.Cells.Item(6, 8).Formula = "=G6/E6"

# Row 7: Osh Park panelized printing quote
$ws.Cells.Item(7, 2).Value = "8/30/2016"
$ws.Cells.Item(7, 5).Value = 12
$ws.Cells.Item(7, 6).Value = 242.5
$ws.Cells.Item(7, 7).Formula = "=F7"
$ws.Cells.Item(7, 8).Formula = "=G7/E7"

# Row 8 (previously row 6, Osh Park prototype) - Cost/Unit formula
$ws.Cells.Item(8, 8).Formula = "=F8/E8"

# Rows 9-11: panelized gerber quote tiers
$ws.Cells.Item(9, 2).Value = "8/29/2016"
$ws.Cells.Item(9, 5).Value = 10
$ws.Cells.Item(9, 6).Value = 2528.49
$ws.Cells.Item(9, 7).Formula = "=F9"
$ws.Cells.Item(9, 8).Formula = "=G9/E9"

$ws.Cells.Item(10, 2).Value = "8/29/2016"
$ws.Cells.Item(10, 5).Value = 20
$ws.Cells.Item(10, 6).Value = 3372.79
$ws.Cells.Item(10, 7).Formula = "=F10"
$ws.Cells.Item(10, 8).Formula = "=G10/E10"

$ws.Cells.Item(11, 2).Value = "8/29/2016"
$ws.Cells.Item(11, 5).Value = 25
$ws.Cells.Item(11, 6).Value = 3826.44
$ws.Cells.Item(11, 7).Formula = "=F11"
$ws.Cells.Item(11, 8).Formula = "=G11/E11"

# ---------------------------------------------------------------
# 5. Text values, written in the exact order new distinct strings
#    were first introduced, so the shared-string table lands in
#    the same order as the canonical workbook.
# ---------------------------------------------------------------

# 276: "4PCB"
$ws.Cells.Item(9, 1).Value = "4PCB"
$ws.Cells.Item(10, 1).Value = "4PCB"
$ws.Cells.Item(11, 1).Value = "4PCB"

# 277: "Assembly+Printing"
$ws.Cells.Item(6, 3).Value = "Assembly+Printing"
$ws.Cells.Item(9, 3).Value = "Assembly+Printing"
$ws.Cells.Item(10, 3).Value = "Assembly+Printing"
$ws.Cells.Item(11, 3).Value = "Assembly+Printing"

# 278: new header in G2
$ws.Cells.Item(2, 7).Value = "Cost Quoted+Extra Parts/Quantities (See formula for details)"

# 279: "10 Day"
$ws.Cells.Item(9, 4).Value = "10 Day"
$ws.Cells.Item(10, 4).Value = "10 Day"
$ws.Cells.Item(11, 4).Value = "10 Day"

# 280: "Includes everything"
$ws.Cells.Item(9, 9).Value = "Includes everything"
$ws.Cells.Item(10, 9).Value = "Includes everything"
$ws.Cells.Item(11, 9).Value = "Includes everything"

# 281: "26 Day"
$ws.Cells.Item(6, 4).Value = "26 Day"

# 282: row 6 note
$ws.Cells.Item(6, 9).Value = "Add approximately `$89.28 for power supply, SD Card holder; 10% Discount for being a student"

# 283: row 7 note
$ws.Cells.Item(7, 9).Value = "Quote is for 3 panels containing 4 pcbs each"

# Remaining text cells that reuse already-existing shared strings
$ws.Cells.Item(6, 1).Value = "Screaming Circuits"
$ws.Cells.Item(7, 1).Value = "Osh Park"
$ws.Cells.Item(7, 3).Value = "Printing"
$ws.Cells.Item(7, 4).Value = "12 Calendar Days"

# ---------------------------------------------------------------
# 6. Column widths
# ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 15.3340625
$ws.Columns.Item(7).ColumnWidth = 52.00203125
$ws.Columns.Item(8).ColumnWidth = 9.50203125
$ws.Columns.Item(9).ColumnWidth = 117.3340625

# ---------------------------------------------------------------
# 7. Selection / active cell
# ---------------------------------------------------------------
$ws.Range("I8").Select()
